$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the rest of row 1 (bold, centered,
# bordered) by copying the style from the neighboring header cell H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells I2 and J2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
